$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-03 Tuesday" "2024-12-04 Wednesday"

Replace-Text "757×6=4542" "631×9=5679"
Replace-Text "692×6=4152" "297×2=594"
Replace-Text "895×5=4475" "454×3=1362"
Replace-Text "621×4=2484" "928×7=6496"
Replace-Text "630×6=3780" "486×5=2430"

Replace-Text "433×9=3897" "564×9=5076"
Replace-Text "611×5=3055" "645×9=5805"
Replace-Text "227×5=1135" "702×6=4212"
Replace-Text "198×4=792" "361×7=2527"
Replace-Text "639×6=3834" "555×4=2220"

Replace-Text "800×2=1600" "256×3=768"
Replace-Text "841×2=1682" "504×2=1008"
Replace-Text "624×6=3744" "401×4=1604"
Replace-Text "184×3=552" "413×5=2065"
Replace-Text "981×9=8829" "903×7=6321"

Replace-Text "813×5=4065" "979×4=3916"
Replace-Text "161×6=966" "657×2=1314"
Replace-Text "793×2=1586" "387×4=1548"
Replace-Text "575×5=2875" "885×5=4425"
Replace-Text "593×3=1779" "972×8=7776"

Replace-Text "401×3=1203" "986×3=2958"
Replace-Text "741×5=3705" "870×9=7830"
Replace-Text "500×5=2500" "193×2=386"
Replace-Text "756×7=5292" "153×5=765"
Replace-Text "932×9=8388" "645×6=3870"
